# Append a new metric data row (row 52) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 52

# Timestamp column (A) is stored as an inline string, not a date value.
$ws.Cells.Item($newRow, 1).Value = "2025-04-29 06:13:31"

# Metric column (B) is a numeric value.
$ws.Cells.Item($newRow, 2).Value = 164

$wb.Save()
